$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 410.36365
$ws.Range("I53").Value = 497.33334
$ws.Range("J53").Value = 350.15384
$ws.Range("K53").Value = 497.33334
$ws.Range("L53").Value = 350.15384
$ws.Range("M53").Value = 139.66666
$ws.Range("N53").Value = -1624.15384
$ws.Range("H82").Value = 3341.8
$ws.Range("I82").Value = 236.33333
$ws.Range("K82").Value = 708.99999
$ws.Range("M82").Value = -302.99999
$ws.Range("H85").Value = 3341.8
$ws.Range("I85").Value = 236.33333
$ws.Range("K85").Value = 708.99999
$ws.Range("M85").Value = 695.00001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7722.253
$ws.Range("I32").Value = 5277.88
$ws.Range("K32").Value = 5277.88
$ws.Range("M32").Value = -4990.88
$ws.Range("H74").Value = 1479.0588
$ws.Range("I74").Value = 1615.909
$ws.Range("J74").Value = 1228.1666
$ws.Range("K74").Value = 1615.909
$ws.Range("L74").Value = 1228.1666
$ws.Range("M74").Value = -741.9090000000001
$ws.Range("N74").Value = -2976.1666
$ws.Range("H77").Value = 1479.0588
$ws.Range("I77").Value = 1615.909
$ws.Range("J77").Value = 1228.1666
$ws.Range("K77").Value = 8079.545
$ws.Range("L77").Value = 6140.833000000001
$ws.Range("M77").Value = -3711.545
$ws.Range("N77").Value = -14876.833
$ws.Range("H102").Value = 2751.6365
$ws.Range("I102").Value = 2026.8
$ws.Range("K102").Value = 2026.8
$ws.Range("M102").Value = -404.8
$ws.Range("H130").Value = 29833.334
$ws.Range("J130").Value = 29833.334
$ws.Range("L130").Value = 29833.334
$ws.Range("N130").Value = -39873.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3535.2727
$ws.Range("I134").Value = 2471.25
$ws.Range("J134").Value = 4812.1
$ws.Range("K134").Value = 7413.75
$ws.Range("L134").Value = 14436.3
$ws.Range("M134").Value = -4878.75
$ws.Range("N134").Value = -19506.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11908276
$ws.Range("I58").Value = 2187.077
$ws.Range("J58").Value = 31255670
$ws.Range("K58").Value = 2187.077
$ws.Range("L58").Value = 31255670
$ws.Range("M58").Value = -1984.077
$ws.Range("N58").Value = -31256076
$ws.Range("H136").Value = 11908276
$ws.Range("I136").Value = 2187.077
$ws.Range("J136").Value = 31255670
$ws.Range("K136").Value = 6561.231000000001
$ws.Range("L136").Value = 93767010
$ws.Range("M136").Value = -4011.231000000001
$ws.Range("N136").Value = -93772110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 291.25
$ws.Range("I4").Value = 66
$ws.Range("J4").Value = 666.6667
$ws.Range("K4").Value = 198
$ws.Range("L4").Value = 2000.0001
$ws.Range("M4").Value = -86
$ws.Range("N4").Value = -2224.0001
$ws.Range("H29").Value = 269.625
$ws.Range("I29").Value = 45
$ws.Range("J29").Value = 301.7143
$ws.Range("K29").Value = 135
$ws.Range("L29").Value = 905.1428999999999
$ws.Range("M29").Value = 142
$ws.Range("N29").Value = -1459.1429
$ws.Range("H55").Value = 1179.6875
$ws.Range("J55").Value = 1537.1
$ws.Range("L55").Value = 4611.299999999999
$ws.Range("N55").Value = -4965.299999999999
$ws.Range("H74").Value = 3273.2856
$ws.Range("I74").Value = 1506.5
$ws.Range("J74").Value = 3980
$ws.Range("K74").Value = 4519.5
$ws.Range("L74").Value = 11940
$ws.Range("M74").Value = -3458.5
$ws.Range("N74").Value = -14062
$ws.Range("H77").Value = 3273.2856
$ws.Range("I77").Value = 1506.5
$ws.Range("J77").Value = 3980
$ws.Range("K77").Value = 13558.5
$ws.Range("L77").Value = 35820
$ws.Range("M77").Value = -8254.5
$ws.Range("N77").Value = -46428
$ws.Range("H87").Value = 9601.154
$ws.Range("I87").Value = 6757.222
$ws.Range("K87").Value = 20271.666
$ws.Range("M87").Value = -19023.666
$ws.Range("H90").Value = 9601.154
$ws.Range("I90").Value = 6757.222
$ws.Range("K90").Value = 60814.998
$ws.Range("M90").Value = -54574.998
$ws.Range("H124").Value = 21386
$ws.Range("I124").Value = 1732.5
$ws.Range("J124").Value = 100000
$ws.Range("K124").Value = 5197.5
$ws.Range("L124").Value = 300000
$ws.Range("M124").Value = -287.5
$ws.Range("N124").Value = -309820
$ws.Range("H138").Value = 2809.3635
$ws.Range("I138").Value = 874.2
$ws.Range("J138").Value = 4422
$ws.Range("K138").Value = 2622.6
$ws.Range("L138").Value = 13266
$ws.Range("M138").Value = 2517.4
$ws.Range("N138").Value = -23546
$ws.Range("H140").Value = 4216.6665
$ws.Range("I140").Value = 1500
$ws.Range("K140").Value = 4500
$ws.Range("M140").Value = 680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4157.0967
$ws.Range("I132").Value = 4363.6
$ws.Range("J132").Value = 3781.6365
$ws.Range("K132").Value = 13090.8
$ws.Range("L132").Value = 11344.9095
$ws.Range("M132").Value = -10560.8
$ws.Range("N132").Value = -16404.9095
$ws.Range("H141").Value = 37944.332
$ws.Range("J141").Value = 37944.332
$ws.Range("N141").Value = -48304.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 4482.5454
$ws.Range("I31").Value = 6786.857
$ws.Range("K31").Value = 6786.857
$ws.Range("M31").Value = -6538.857
$ws.Range("H132").Value = 2180.9736
$ws.Range("I132").Value = 1282.2727
$ws.Range("J132").Value = 3416.6875
$ws.Range("K132").Value = 3846.8181
$ws.Range("L132").Value = 10250.0625
$ws.Range("M132").Value = -1316.8181
$ws.Range("N132").Value = -15310.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 13475.765
$ws.Range("I58").Value = 6056.6665
$ws.Range("J58").Value = 13939.458
$ws.Range("K58").Value = 6056.6665
$ws.Range("L58").Value = 13939.458
$ws.Range("M58").Value = -5748.6665
$ws.Range("N58").Value = -14555.458
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("H122").Value = 627208.1
$ws.Range("I122").Value = 668355.3
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 2005065.9
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -2002615.9
$ws.Range("N122").Value = -34900
$ws.Range("H136").Value = 1317.2667
$ws.Range("I136").Value = 596.5333000000001
$ws.Range("J136").Value = 2038
$ws.Range("K136").Value = 1789.5999
$ws.Range("L136").Value = 6114
$ws.Range("M136").Value = 760.4000999999998
$ws.Range("N136").Value = -11214
$ws.Range("N112").ClearContents()

